# Apply updated crypto price/volume data to the worksheet.
# Source: scheduled GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.249.88"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "2.579.15"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  +2.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.25"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "2.585.40"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("E11").Value = "  +2.86%  "
$ws.Range("E12").Value = "  +11.16%  "
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("D14").Value = "3.033.98"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").Value = "59.248.68"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.48"
$ws.Range("E16").Value = "  +6.93%  "
$ws.Range("E17").Value = "  +3.58%  "
$ws.Range("D18").Value = "2.584.26"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "339.53"
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.26"
$ws.Range("E21").Value = "  +1.44%  "
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.58"
$ws.Range("E24").Value = "  -2.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.461"
$ws.Range("E25").Value = "  +7.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.993"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.28"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("E29").Value = "  +2.83%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.06"
$ws.Range("E32").Value = "  +1.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.05"
$ws.Range("E33").Value = "  +2.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.03"
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("E35").Value = "  +1.74%  "
$ws.Range("E36").Value = "  +2.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.877"
$ws.Range("E37").Value = "  -3.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.871"
$ws.Range("E38").Value = "  -4.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.29"
$ws.Range("E39").Value = "  +0.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "296.78"
$ws.Range("E41").Value = "  +4.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.67"
$ws.Range("E42").Value = "  +2.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0977"
$ws.Range("E44").Value = "  +2.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "129.54"
$ws.Range("E45").Value = "  +10.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.594"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0537"
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.21"
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.64"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("D51").Value = "1.956.33"
$ws.Range("E51").Value = "  +0.66%  "
